# Auto-generated Excel COM-interop script
# Applies cell value updates across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2749.75
$ws.Range("J48").Value = 2749.75
$ws.Range("L48").Value = 8249.25
$ws.Range("N48").Value = -8833.25
$ws.Range("H56").Value = 2749.75
$ws.Range("J56").Value = 2749.75
$ws.Range("L56").Value = 8249.25
$ws.Range("N56").Value = -9317.25
$ws.Range("H112").Value = 3573084.2
$ws.Range("I112").Value = 7143714.5
$ws.Range("J112").Value = 2454.2856
$ws.Range("K112").Value = 21431143.5
$ws.Range("L112").Value = 7362.8568
$ws.Range("M112").Value = -21430035.5
$ws.Range("N112").Value = -9578.856800000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5680.577
$ws.Range("I32").Value = 5307.8
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 5307.8
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -5020.8
$ws.Range("N32").Value = -15574
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H61").Value = 5555.273
$ws.Range("I61").Value = 5110.8
$ws.Range("K61").Value = 5110.8
$ws.Range("M61").Value = -4898.8
$ws.Range("H63").Value = 1867.6842
$ws.Range("I63").Value = 1168.1538
$ws.Range("K63").Value = 1168.1538
$ws.Range("M63").Value = -482.1538
$ws.Range("H66").Value = 1867.6842
$ws.Range("I66").Value = 1168.1538
$ws.Range("K66").Value = 5840.769
$ws.Range("M66").Value = -2408.769
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 1821.3334
$ws.Range("I132").Value = 1855.75
$ws.Range("J132").Value = 1752.5
$ws.Range("K132").Value = 5567.25
$ws.Range("L132").Value = 5257.5
$ws.Range("M132").Value = -3037.25
$ws.Range("N132").Value = -10317.5
$ws.Range("H136").Value = 5555.273
$ws.Range("I136").Value = 5110.8
$ws.Range("K136").Value = 15332.4
$ws.Range("M136").Value = -12782.4

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 255.5
$ws.Range("I12").Value = 142.71428
$ws.Range("J12").Value = 518.6667
$ws.Range("K12").Value = 142.71428
$ws.Range("L12").Value = 518.6667
$ws.Range("M12").Value = 25.28572
$ws.Range("N12").Value = -854.6667
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 200
$ws.Range("K23").Value = 200
$ws.Range("M23").Value = 83
$ws.Range("H25").Value = 549
$ws.Range("I25").Value = 549
$ws.Range("K25").Value = 549
$ws.Range("M25").Value = -314
$ws.Range("H105").Value = 3133.3333
$ws.Range("I105").Value = 2950
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 2950
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -1203
$ws.Range("N105").Value = -6994
$ws.Range("H134").Value = 8848.762000000001
$ws.Range("I134").Value = 8096
$ws.Range("K134").Value = 24288
$ws.Range("M134").Value = -21753

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 35010.9
$ws.Range("I4").Value = 35010.9
$ws.Range("K4").Value = 35010.9
$ws.Range("M4").Value = -34898.9
$ws.Range("H6").Value = 353.45456
$ws.Range("I6").Value = 385.8
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 385.8
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = -272.8
$ws.Range("N6").Value = -256
$ws.Range("H31").Value = 5540.923
$ws.Range("I31").Value = 4093.1428
$ws.Range("K31").Value = 4093.1428
$ws.Range("M31").Value = -3798.1428
$ws.Range("H34").Value = 5540.923
$ws.Range("I34").Value = 4093.1428
$ws.Range("K34").Value = 4093.1428
$ws.Range("M34").Value = -3891.1428
$ws.Range("H58").Value = 2231.8572
$ws.Range("J58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10406
$ws.Range("H136").Value = 2231.8572
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -35100

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 305.93332
$ws.Range("I7").Value = 139.33333
$ws.Range("J7").Value = 417
$ws.Range("K7").Value = 417.99999
$ws.Range("L7").Value = 1251
$ws.Range("M7").Value = -305.99999
$ws.Range("N7").Value = -1475
$ws.Range("H62").Value = 9500
$ws.Range("I62").Value = 9500
$ws.Range("K62").Value = 28500
$ws.Range("M62").Value = -27814
$ws.Range("H65").Value = 9500
$ws.Range("I65").Value = 9500
$ws.Range("K65").Value = 85500
$ws.Range("M65").Value = -82068

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 890
$ws.Range("I4").Value = 2600
$ws.Range("J4").Value = 320
$ws.Range("K4").Value = 2600
$ws.Range("L4").Value = 320
$ws.Range("M4").Value = -2488
$ws.Range("N4").Value = -544
$ws.Range("H5").Value = 1300.1
$ws.Range("I5").Value = 1300.1
$ws.Range("K5").Value = 1300.1
$ws.Range("M5").Value = -1188.1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H13").Value = 399.77777
$ws.Range("I13").Value = 240
$ws.Range("J13").Value = 599.5
$ws.Range("K13").Value = 240
$ws.Range("L13").Value = 599.5
$ws.Range("M13").Value = -101
$ws.Range("N13").Value = -877.5
$ws.Range("H17").Value = 736.3333
$ws.Range("J17").Value = 736.3333
$ws.Range("L17").Value = 736.3333
$ws.Range("N17").Value = -1072.3333
$ws.Range("H19").Value = 815.1111
$ws.Range("I19").Value = 600.8333
$ws.Range("J19").Value = 1243.6666
$ws.Range("K19").Value = 600.8333
$ws.Range("L19").Value = 1243.6666
$ws.Range("M19").Value = -312.8333
$ws.Range("N19").Value = -1819.6666
$ws.Range("H22").Value = 1950
$ws.Range("J22").Value = 1950
$ws.Range("L22").Value = 1950
$ws.Range("N22").Value = -3008
$ws.Range("H23").Value = 818
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 818
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 818
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1264
$ws.Range("H24").Value = 37855.9
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 37855.9
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 37855.9
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -38201.9
$ws.Range("H25").Value = 3000
$ws.Range("J25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -4058
$ws.Range("H29").Value = 6995
$ws.Range("J29").Value = 6995
$ws.Range("L29").Value = 6995
$ws.Range("N29").Value = -7575
$ws.Range("H113").Value = 4333.3335
$ws.Range("J113").Value = 8000
$ws.Range("L113").Value = 8000
$ws.Range("N113").Value = -12340
$ws.Range("H132").Value = 2862.25
$ws.Range("I132").Value = 2522.7693
$ws.Range("K132").Value = 7568.3079
$ws.Range("M132").Value = -5038.3079

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6298.5
$ws.Range("I132").Value = 4583.3335
$ws.Range("J132").Value = 7584.875
$ws.Range("K132").Value = 13750.0005
$ws.Range("L132").Value = 22754.625
$ws.Range("M132").Value = -11220.0005
$ws.Range("N132").Value = -27814.625
$ws.Range("H136").Value = 3089.889
$ws.Range("I136").Value = 3089.889
$ws.Range("K136").Value = 9269.667000000001
$ws.Range("M136").Value = -6719.667000000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 877.44446
$ws.Range("I122").Value = 874
$ws.Range("K122").Value = 2622
$ws.Range("M122").Value = -172

Write-Host "Applied 210 cell updates across 8 sheets."